$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily coinranking.com data refresh (2023-01-15 18:40 UTC GitHub Actions run):
# updated Price (column D) and Volume(1h) (column E) figures for each coin row.
# NumberFormat is forced to Text ("@") before each write and cleared afterwards so
# that the numeric-/percent-looking strings are stored as plain text, exactly like
# the rest of the sheet (which already stores these columns as text).
$cells = @(
    @{Addr = "D2"; Val = "301.31"}
    @{Addr = "E2"; Val = "-0.64%"}
    @{Addr = "D3"; Val = "31.43"}
    @{Addr = "E3"; Val = "-1.68%"}
    @{Addr = "D4"; Val = "5.161"}
    @{Addr = "E4"; Val = "-1.86%"}
    @{Addr = "E5"; Val = "-1.10%"}
    @{Addr = "D6"; Val = "2.614"}
    @{Addr = "E6"; Val = "72.11%"}
    @{Addr = "D7"; Val = "7.911"}
    @{Addr = "E7"; Val = "0.72%"}
    @{Addr = "D8"; Val = "3.753"}
    @{Addr = "E8"; Val = "-1.23%"}
    @{Addr = "D9"; Val = "0.9202"}
    @{Addr = "E9"; Val = "0.12%"}
    @{Addr = "D10"; Val = "0.1738"}
    @{Addr = "E10"; Val = "3.07%"}
    @{Addr = "D11"; Val = "0.07490"}
    @{Addr = "E11"; Val = "-6.76%"}
    @{Addr = "D12"; Val = "0.08138"}
    @{Addr = "E12"; Val = "2.49%"}
    @{Addr = "D13"; Val = "0.03047"}
    @{Addr = "E13"; Val = "0.56%"}
    @{Addr = "D14"; Val = "0.09916"}
    @{Addr = "E14"; Val = "0.25%"}
    @{Addr = "D15"; Val = "0.001509"}
    @{Addr = "E15"; Val = "0.95%"}
    @{Addr = "D16"; Val = "0.006102"}
    @{Addr = "E16"; Val = "-1.94%"}
    @{Addr = "D17"; Val = "3.448"}
    @{Addr = "E17"; Val = "-0.68%"}
    @{Addr = "D18"; Val = "2.230"}
    @{Addr = "E18"; Val = "-0.02%"}
    @{Addr = "E19"; Val = "-1.06%"}
    @{Addr = "D20"; Val = "0.1350"}
    @{Addr = "E20"; Val = "1.48%"}
    @{Addr = "D21"; Val = "4.660"}
    @{Addr = "E21"; Val = "3.91%"}
    @{Addr = "D22"; Val = "0.04654"}
    @{Addr = "E22"; Val = "0.78%"}
    @{Addr = "D23"; Val = "0.1570"}
    @{Addr = "E23"; Val = "-3.24%"}
    @{Addr = "E24"; Val = "0.52%"}
    @{Addr = "D25"; Val = "0.004475"}
    @{Addr = "E25"; Val = "0.61%"}
    @{Addr = "E26"; Val = "-7.02%"}
    @{Addr = "E27"; Val = "7.19%"}
    @{Addr = "D39"; Val = "0.01723"}
    @{Addr = "E39"; Val = "-1.29%"}
    @{Addr = "D40"; Val = "0.04522"}
    @{Addr = "E40"; Val = "0.74%"}
    @{Addr = "D41"; Val = "0.007194"}
    @{Addr = "E41"; Val = "0.31%"}
    @{Addr = "D42"; Val = "0.1343"}
    @{Addr = "E42"; Val = "-0.32%"}
    @{Addr = "D43"; Val = "0.002213"}
    @{Addr = "E43"; Val = "0.13%"}
    @{Addr = "D44"; Val = "0.01091"}
    @{Addr = "E44"; Val = "-14.51%"}
    @{Addr = "D45"; Val = "0.00006285"}
    @{Addr = "E45"; Val = "2.53%"}
    @{Addr = "D46"; Val = "0.01001"}
    @{Addr = "E46"; Val = "-22.99%"}
    @{Addr = "D47"; Val = "1.928"}
    @{Addr = "E47"; Val = "3.20%"}
)

foreach ($item in $cells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.ClearFormats()
}

